$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 43 and 44: the two fixtures (Espanyol-Amorebieta / Zaragoza-Eldense)
#    had their match-data columns (F..V) swapped between the two rows.
#    Columns A..E (index/pais/torneio/temporada/data_partida) stay as-is.
# ---------------------------------------------------------------------------
$row43 = @("Zaragoza", 2, "Eldense", 0, 1.68, "29/08/2023 08:42", 1.83, "03/09/2023 18:16", `
           3.54, "29/08/2023 08:42", 3.28, "03/09/2023 18:24", 6.01, "29/08/2023 08:42", `
           5.59, "03/09/2023 18:16", "https://www.betexplorer.com/football/spain/laliga2/zaragoza-eldense/dxXyUnaA/")

$row44 = @("Espanyol", 3, "Amorebieta", 2, 1.43, "28/08/2023 23:42", 1.54, "03/09/2023 18:25", `
           4.54, "28/08/2023 23:42", 3.99, "03/09/2023 18:28", 7.98, "28/08/2023 23:42", `
           7.65, "03/09/2023 18:28", "https://www.betexplorer.com/football/spain/laliga2/espanyol-amorebieta/6mvPWlFi/")

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "43").Value2 = $row43[$i]
    $ws.Range($cols[$i] + "44").Value2 = $row44[$i]
}

# ---------------------------------------------------------------------------
# 2) Rows 72 and 73: the two fixtures (Alcorcon-Huesca / Racing Santander-
#    Albacete) also had their match-data columns (F..V) swapped.
# ---------------------------------------------------------------------------
$row72 = @("Racing Santander", 2, "Albacete", 1, 2.56, "16/09/2023 20:12", 2.66, "23/09/2023 20:57", `
           3.12, "16/09/2023 20:12", 3, "23/09/2023 20:57", 3.12, "16/09/2023 20:12", `
           3.15, "23/09/2023 20:57", "https://www.betexplorer.com/football/spain/laliga2/racing-santander-albacete/l6CO8rS8/")

$row73 = @("Alcorcon", 0, "Huesca", 2, 2.4, "18/09/2023 11:42", 2.31, "23/09/2023 20:55", `
           2.95, "18/09/2023 11:42", 2.91, "23/09/2023 20:56", 3.54, "18/09/2023 11:42", `
           3.98, "23/09/2023 20:56", "https://www.betexplorer.com/football/spain/laliga2/alcorcon-huesca/2aDQmsjq/")

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "72").Value2 = $row72[$i]
    $ws.Range($cols[$i] + "73").Value2 = $row73[$i]
}

# ---------------------------------------------------------------------------
# 3) Two brand-new fixture rows were appended at the bottom (77 and 78).
#    Copy the formatting of the last existing row (76) down into the new
#    rows so the styles (bold/bordered index column, date-formatted data
#    column) match, then overwrite the values.
# ---------------------------------------------------------------------------
$ws.Range("A76:V76").Copy($ws.Range("A77:V77"))
$ws.Range("A76:V76").Copy($ws.Range("A78:V78"))

$ws.Range("A77").Value2 = 76
$ws.Range("B77").Value2 = "spain"
$ws.Range("C77").Value2 = "laliga2"
$ws.Range("D77").Value2 = "2023-2024"
$ws.Range("E77").Value2 = 45194.875
$ws.Range("F77").Value2 = "Ferrol"
$ws.Range("G77").Value2 = 1
$ws.Range("H77").Value2 = "Zaragoza"
$ws.Range("I77").Value2 = 0
$ws.Range("J77").Value2 = 2.69
$ws.Range("K77").Value2 = "18/09/2023 11:42"
$ws.Range("L77").Value2 = 2.26
$ws.Range("M77").Value2 = "25/09/2023 20:59"
$ws.Range("N77").Value2 = 3.01
$ws.Range("O77").Value2 = "18/09/2023 11:42"
$ws.Range("P77").Value2 = 3.06
$ws.Range("Q77").Value2 = "25/09/2023 20:59"
$ws.Range("R77").Value2 = 3
$ws.Range("S77").Value2 = "18/09/2023 11:42"
$ws.Range("T77").Value2 = 3.85
$ws.Range("U77").Value2 = "25/09/2023 20:59"
$ws.Range("V77").Value2 = "https://www.betexplorer.com/football/spain/laliga2/ferrol-zaragoza/YcUyouz2/"

$ws.Range("A78").Value2 = 77
$ws.Range("B78").Value2 = "spain"
$ws.Range("C78").Value2 = "laliga2"
$ws.Range("D78").Value2 = "2023-2024"
$ws.Range("E78").Value2 = 45194.875
$ws.Range("F78").Value2 = "Tenerife"
$ws.Range("G78").Value2 = 1
$ws.Range("H78").Value2 = "Espanyol"
$ws.Range("I78").Value2 = 0
$ws.Range("J78").Value2 = 2.62
$ws.Range("K78").Value2 = "24/09/2023 16:13"
$ws.Range("L78").Value2 = 2.73
$ws.Range("M78").Value2 = "25/09/2023 20:20"
$ws.Range("N78").Value2 = 3.07
$ws.Range("O78").Value2 = "24/09/2023 16:13"
$ws.Range("P78").Value2 = 2.99
$ws.Range("Q78").Value2 = "25/09/2023 20:20"
$ws.Range("R78").Value2 = 3.02
$ws.Range("S78").Value2 = "24/09/2023 16:13"
$ws.Range("T78").Value2 = 3.06
$ws.Range("U78").Value2 = "25/09/2023 20:20"
$ws.Range("V78").Value2 = "https://www.betexplorer.com/football/spain/laliga2/tenerife-espanyol/xUT0kfRm/"

# ---------------------------------------------------------------------------
# 4) Refresh the sheet dimension to cover the two newly-added rows.
# ---------------------------------------------------------------------------
$ws.UsedRange | Out-Null
